$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 194444910
$ws.Range("J43").Value = 55556092
$ws.Range("L43").Value = 55556092
$ws.Range("N43").Value = -55556230

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1857.2549
$ws.Range("I58").Value = 338.36365
$ws.Range("J58").Value = 2274.95
$ws.Range("K58").Value = 1015.09095
$ws.Range("L58").Value = 6824.849999999999
$ws.Range("M58").Value = -865.09095
$ws.Range("N58").Value = -7124.849999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2281.4285
$ws.Range("I64").Value = 2219.5
$ws.Range("J64").Value = 2479.6
$ws.Range("K64").Value = 2219.5
$ws.Range("L64").Value = 2479.6
$ws.Range("M64").Value = -1971.5
$ws.Range("N64").Value = -2975.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2281.4285
$ws.Range("I67").Value = 2219.5
$ws.Range("J67").Value = 2479.6
$ws.Range("K67").Value = 2219.5
$ws.Range("L67").Value = 2479.6
$ws.Range("M67").Value = -1361.5
$ws.Range("N67").Value = -4195.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 50531.285
$ws.Range("I76").Value = 52907.65
$ws.Range("J76").Value = 3004
$ws.Range("K76").Value = 52907.65
$ws.Range("L76").Value = 3004
$ws.Range("M76").Value = -52592.65
$ws.Range("N76").Value = -3634

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 50531.285
$ws.Range("I79").Value = 52907.65
$ws.Range("J79").Value = 3004
$ws.Range("K79").Value = 52907.65
$ws.Range("L79").Value = 3004
$ws.Range("M79").Value = -51815.65
$ws.Range("N79").Value = -5188

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 21784740
$ws.Range("I106").Value = 54526.316
$ws.Range("K106").Value = 54526.316
$ws.Range("M106").Value = -53895.316

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1906342.2
$ws.Range("I132").Value = 1932036.1
$ws.Range("K132").Value = 5796108.300000001
$ws.Range("M132").Value = -5793578.300000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 25001182
$ws.Range("I137").Value = 920.55884
$ws.Range("J137").Value = 166669330
$ws.Range("K137").Value = 2761.67652
$ws.Range("L137").Value = 500007990
$ws.Range("M137").Value = -211.67652
$ws.Range("N137").Value = -500013090

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3310.8923
$ws.Range("I138").Value = 2889.524
$ws.Range("J138").Value = 3512
$ws.Range("K138").Value = 8668.572
$ws.Range("L138").Value = 10536
$ws.Range("M138").Value = -3528.572
$ws.Range("N138").Value = -20816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30467.6
$ws.Range("I32").Value = 33259.85
$ws.Range("J32").Value = 26000
$ws.Range("K32").Value = 33259.85
$ws.Range("L32").Value = 26000
$ws.Range("M32").Value = -32972.85
$ws.Range("N32").Value = -26574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1418.4386
$ws.Range("I61").Value = 1164.6171
$ws.Range("J61").Value = 2611.4
$ws.Range("K61").Value = 1164.6171
$ws.Range("L61").Value = 2611.4
$ws.Range("M61").Value = -952.6170999999999
$ws.Range("N61").Value = -3035.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3093.625
$ws.Range("I63").Value = 2476.182
$ws.Range("J63").Value = 4452
$ws.Range("K63").Value = 2476.182
$ws.Range("L63").Value = 4452
$ws.Range("M63").Value = -1790.182
$ws.Range("N63").Value = -5824

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3093.625
$ws.Range("I66").Value = 2476.182
$ws.Range("J66").Value = 4452
$ws.Range("K66").Value = 12380.91
$ws.Range("L66").Value = 22260
$ws.Range("M66").Value = -8948.91
$ws.Range("N66").Value = -29124

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1188.5834
$ws.Range("I74").Value = 1136
$ws.Range("J74").Value = 1346.3334
$ws.Range("K74").Value = 1136
$ws.Range("L74").Value = 1346.3334
$ws.Range("M74").Value = -262
$ws.Range("N74").Value = -3094.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 25247
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 25247
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 25247
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -25923

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1188.5834
$ws.Range("I77").Value = 1136
$ws.Range("J77").Value = 1346.3334
$ws.Range("K77").Value = 5680
$ws.Range("L77").Value = 6731.666999999999
$ws.Range("M77").Value = -1312
$ws.Range("N77").Value = -15467.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 25247
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 25247
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 25247
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -27587

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6312.745
$ws.Range("I132").Value = 7271.5405
$ws.Range("K132").Value = 21814.6215
$ws.Range("M132").Value = -19284.6215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1418.4386
$ws.Range("I136").Value = 1164.6171
$ws.Range("J136").Value = 2611.4
$ws.Range("K136").Value = 3493.8513
$ws.Range("L136").Value = 7834.200000000001
$ws.Range("M136").Value = -943.8512999999998
$ws.Range("N136").Value = -12934.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 49633.332
$ws.Range("J130").Value = 49633.332
$ws.Range("L130").Value = 49633.332
$ws.Range("N130").Value = -59673.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4234.4814
$ws.Range("I134").Value = 5744.1724
$ws.Range("J134").Value = 2483.24
$ws.Range("K134").Value = 17232.5172
$ws.Range("L134").Value = 7449.719999999999
$ws.Range("M134").Value = -14697.5172
$ws.Range("N134").Value = -12519.72

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4500.0815
$ws.Range("I31").Value = 2530
$ws.Range("J31").Value = 6242.846
$ws.Range("K31").Value = 2530
$ws.Range("L31").Value = 6242.846
$ws.Range("M31").Value = -2235
$ws.Range("N31").Value = -6832.846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4500.0815
$ws.Range("I34").Value = 2530
$ws.Range("J34").Value = 6242.846
$ws.Range("K34").Value = 2530
$ws.Range("L34").Value = 6242.846
$ws.Range("M34").Value = -2328
$ws.Range("N34").Value = -6646.846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5953747
$ws.Range("I58").Value = 1241.9762
$ws.Range("J58").Value = 23811262
$ws.Range("K58").Value = 1241.9762
$ws.Range("L58").Value = 23811262
$ws.Range("M58").Value = -1038.9762
$ws.Range("N58").Value = -23811668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2710.5103
$ws.Range("I134").Value = 2710.5103
$ws.Range("K134").Value = 8131.5309
$ws.Range("M134").Value = -5596.5309

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5953747
$ws.Range("I136").Value = 1241.9762
$ws.Range("J136").Value = 23811262
$ws.Range("K136").Value = 3725.9286
$ws.Range("L136").Value = 71433786
$ws.Range("M136").Value = -1175.9286
$ws.Range("N136").Value = -71438886

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 723.5714
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 1666
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 4998
$ws.Range("M4").Value = -488
$ws.Range("N4").Value = -5222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2542272.8
$ws.Range("I131").Value = 33886.668
$ws.Range("J131").Value = 2686987.2
$ws.Range("K131").Value = 101660.004
$ws.Range("L131").Value = 8060961.600000001
$ws.Range("M131").Value = -96620.00399999999
$ws.Range("N131").Value = -8071041.600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 18520374
$ws.Range("I122").Value = 25642588
$ws.Range("K122").Value = 76927764
$ws.Range("M122").Value = -76925314

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6804.75
$ws.Range("I132").Value = 7552.737
$ws.Range("J132").Value = 3962.4
$ws.Range("K132").Value = 22658.211
$ws.Range("L132").Value = 11887.2
$ws.Range("M132").Value = -20128.211
$ws.Range("N132").Value = -16947.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4369.091
$ws.Range("I136").Value = 4816.6763
$ws.Range("J136").Value = 2847.3
$ws.Range("K136").Value = 14450.0289
$ws.Range("L136").Value = 8541.900000000001
$ws.Range("M136").Value = -11900.0289
$ws.Range("N136").Value = -13641.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1362.2
$ws.Range("I136").Value = 1442.6364
$ws.Range("J136").Value = 1141
$ws.Range("K136").Value = 4327.9092
$ws.Range("L136").Value = 3423
$ws.Range("M136").Value = -1777.9092
$ws.Range("N136").Value = -8523
